$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 1. E3 gets a new value "XXX" (was blank)
$ws.Range("E3").Value = "XXX"

# 2. E1 comment/help text gets an extra line describing the formula used
$cell = $ws.Range("E1")

$line0 = "0. Must be 'XXX' or 2 numbers`n"
$line1 = "1. Custom validation as columns C and D`n2. Using formula `n=OR(E4=`"XXX`";AND(LEN(E4)= 2; ISNUMBER(E4)))"

$cell.Value = $line0 + $line1

# Re-apply the original (non-bold) formatting to the second run, since
# setting .Value resets the whole cell to its default (bold) font.
$run2 = $cell.Characters($line0.Length + 1, $line1.Length)
$run2.Font.Bold = $false
$run2.Font.Size = 14
$run2.Font.Color = 0
$run2.Font.Name = "Calibri"

# 3. Move the active selection from C17 to E11 (Sheet1's own remembered
#    selection), while leaving the AUX tab as the workbook's active tab,
#    matching the original file.
$ws.Range("E11").Select()
$wb.Worksheets.Item("AUX").Activate()
